$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New "MILESTONE" lookup block (mirrors the existing ROUTE/ROLE tables) ---
# Row 46: header row (ROLE / ROUTE / CREATE 1 / DELETE 1 / UPDATE 1 / READ 1 /
#          BULK CREATE / BULK DELETE / BULK UPDATE / BULK READ)
$ws.Range("A46").Value = "ROLE"
$ws.Range("B46").Value = "ROUTE"
$ws.Range("C46").Value = "CREATE 1"
$ws.Range("D46").Value = "DELETE 1"
$ws.Range("E46").Value = "UPDATE 1"
$ws.Range("F46").Value = "READ 1"
$ws.Range("G46").Value = "BULK CREATE"
$ws.Range("H46").Value = "BULK DELETE"
$ws.Range("I46").Value = "BULK UPDATE"
$ws.Range("J46").Value = "BULK READ"

# Rows 47-49: STUDENT / TEACHER / ADMIN, with the new MILESTONE route name in
# the merged B47:B49 cell.
$ws.Range("A47").Value = "STUDENT"
$ws.Range("B47").Value = "MILESTONE"
$ws.Range("A48").Value = "TEACHER"
$ws.Range("A49").Value = "ADMIN"

# Merge the role-name cell across the 3 data rows, like every other block,
# before copying formats over it below (Merge can otherwise perturb the
# border of the cells it spans).
$ws.Range("B47:B49").Merge()

# --- Formatting: copy the cell styles from matching cells elsewhere on the
# sheet that already use the required style so no new style entries are
# created (matches the target workbook, whose styles.xml is untouched). ---
function CopyFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial(-4122)  # xlPasteFormats
}

CopyFormat "A1" "A46"
CopyFormat "B1" "B46:F46"
CopyFormat "G1" "G46:J46"

CopyFormat "A32" "A47:A49"
CopyFormat "B32" "B47:B49"
CopyFormat "E27" "C47:F48"
CopyFormat "E27" "C49:F49"
CopyFormat "G37" "G47:H48"
CopyFormat "G7" "I47:I49"
CopyFormat "D2" "J47:J48"
CopyFormat "G38" "G49:H49"
CopyFormat "C2" "J49"

$excel.CutCopyMode = $false

# --- Restore the selection to where the author last clicked ---
$ws.Range("L16").Select()
